$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (type/constraint row) for columns whose value changes ---
$ws.Range("J2").Value = "#integer,  unit:nm"
$ws.Range("M2").Value = "#float"
$ws.Range("N2").Value = "#float"
$ws.Range("O2").Value = "#float,  unit:l/mol/cm"
$ws.Range("P2").Value = "#float"
$ws.Range("I2").Value = "#float,  unit:mlormg"

# --- Add new row 3 (description row) ---
# Columns A-G get the field descriptions; H-P have no description text
# (they stay blank, same as the source row only carries data through G).
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
